$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'69.830.31"
$ws.Range("E2").Value = "  -0.18%  "

# Row 3
$ws.Range("D3").Value = "'3.518.05"
$ws.Range("E3").Value = "  +0.10%  "

# Row 4
$ws.Range("E4").Value = "  -0.06%  "

# Row 5
$ws.Range("D5").Value = "'602.78"
$ws.Range("E5").Value = "  -0.92%  "

# Row 6
$ws.Range("D6").Value = "'195.80"
$ws.Range("E6").Value = "  +1.46%  "

# Row 7
$ws.Range("D7").Value = "'0.623"
$ws.Range("E7").Value = "  -1.00%  "

# Row 8
$ws.Range("E8").Value = "  +0.03%  "

# Row 9
$ws.Range("D9").Value = "'0.203"
$ws.Range("E9").Value = "  -5.15%  "

# Row 10
$ws.Range("D10").Value = "'0.645"
$ws.Range("E10").Value = "  -3.06%  "

# Row 11
$ws.Range("D11").Value = "'53.23"
$ws.Range("E11").Value = "  -0.88%  "

# Row 12
$ws.Range("D12").Value = "'0.0000302"
$ws.Range("E12").Value = "  -2.47%  "

# Row 13
$ws.Range("D13").Value = "'9.48"
$ws.Range("E13").Value = "  -1.78%  "

# Row 14
$ws.Range("D14").Value = "'4.076.61"
$ws.Range("E14").Value = "  -0.11%  "

# Row 15
$ws.Range("D15").Value = "'602.54"
$ws.Range("E15").Value = "  -2.50%  "

# Row 16
$ws.Range("D16").Value = "'69.954.15"
$ws.Range("E16").Value = "  -0.16%  "

# Row 17
$ws.Range("D17").Value = "'12.71"
$ws.Range("E17").Value = "  -0.48%  "

# Row 18
$ws.Range("D18").Value = "'19.02"
$ws.Range("E18").Value = "  +0.39%  "

# Row 19
$ws.Range("B19").Value = "WrappedEther"
$ws.Range("C19").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D19").Value = "'3.514.65"
$ws.Range("E19").Value = "  -0.04%  "

# Row 20
$ws.Range("B20").Value = "TRON"
$ws.Range("C20").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D20").Value = "'0.123"
$ws.Range("E20").Value = "  +1.88%  "

# Row 21
$ws.Range("D21").Value = "'0.982"
$ws.Range("E21").Value = "  -1.27%  "

# Row 22
$ws.Range("D22").Value = "'18.01"
$ws.Range("E22").Value = "  +1.46%  "

# Row 23
$ws.Range("E23").Value = "  +3.04%  "

# Row 24
$ws.Range("D24").Value = "'103.39"
$ws.Range("E24").Value = "  -2.68%  "

# Row 25
$ws.Range("E25").Value = "  -1.34%  "

# Row 26
$ws.Range("E26").Value = "  -0.01%  "

# Row 27
$ws.Range("D27").Value = "'10.78"
$ws.Range("E27").Value = "  -2.10%  "

# Row 28
$ws.Range("D28").Value = "'9.54"
$ws.Range("E28").Value = "  -3.36%  "

# Row 29
$ws.Range("D29").Value = "'33.20"
$ws.Range("E29").Value = "  -2.89%  "

# Row 30
$ws.Range("D30").Value = "'4.41"
$ws.Range("E30").Value = "  +6.96%  "

# Row 31
$ws.Range("D31").Value = "'7.08"
$ws.Range("E31").Value = "  -0.11%  "

# Row 32
$ws.Range("D32").Value = "'12.32"
$ws.Range("E32").Value = "  -1.32%  "

# Row 33
$ws.Range("D33").Value = "'0.114"
$ws.Range("E33").Value = "  -1.22%  "

# Row 34
$ws.Range("D34").Value = "'63.16"
$ws.Range("E34").Value = "  -1.38%  "

# Row 35
$ws.Range("D35").Value = "'3.19"
$ws.Range("E35").Value = "  +2.45%  "

# Row 36
$ws.Range("D36").Value = "'3.755.61"
$ws.Range("E36").Value = "  +1.84%  "

# Row 37
$ws.Range("D37").Value = "'0.0₃0820"
$ws.Range("E37").Value = "  +4.65%  "

# Row 38
$ws.Range("E38").Value = "  -0.02%  "

# Row 39
$ws.Range("E39").Value = "  -0.15%  "

# Row 40
$ws.Range("B40").Value = "TheGraph"
$ws.Range("C40").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D40").Value = "'0.389"
$ws.Range("E40").Value = "  -1.01%  "

# Row 41
$ws.Range("B41").Value = "Bittensor"
$ws.Range("C41").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D41").Value = "'503.15"
$ws.Range("E41").Value = "  -2.30%  "

# Row 42
$ws.Range("D42").Value = "'36.27"
$ws.Range("E42").Value = "  -1.59%  "

# Row 43
$ws.Range("D43").Value = "'0.134"
$ws.Range("E43").Value = "  -3.40%  "

# Row 44
$ws.Range("E44").Value = "  -3.58%  "

# Row 45
$ws.Range("B45").Value = "ThetaToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D45").Value = "'2.82"
$ws.Range("E45").Value = "  -2.89%  "

# Row 46
$ws.Range("B46").Value = "Stellar"
$ws.Range("C46").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D46").Value = "'0.139"
$ws.Range("E46").Value = "  -1.81%  "

# Row 47
$ws.Range("D47").Value = "'3.29"
$ws.Range("E47").Value = "  -0.75%  "

# Row 48
$ws.Range("E48").Value = "  -0.05%  "

# Row 49
$ws.Range("D49").Value = "'8.47"
$ws.Range("E49").Value = "  -3.33%  "

# Row 50
$ws.Range("D50").Value = "'0.000250"
$ws.Range("E50").Value = "  +4.35%  "

# Row 51
$ws.Range("E51").Value = "  +4.92%  "
